$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Metadata sheet: URL / Version / Date / Publisher updates
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/care-gap-compliance-frequency"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# ---------------------------------------------------------------------------
# Elements sheet: the generated snapshot table shuffled which cached text
# blocks line up with which element rows. Row 2 (the root "Extension" row)
# no longer carries the ele-1/ext-1 constraint text (it moved down to row 4,
# "Extension.extension"); row 3 became "Extension.id" and row 4 became
# "Extension.extension" (were previously off-by-one).
# ---------------------------------------------------------------------------
$els = $wb.Worksheets.Item("Elements")

# Row 2 ("Extension"): Constraint(s) cell is now blank.
$els.Range("AI2").Value = ""

# Row 3: now describes "Extension.id"
$els.Range("A3").Value  = "Extension.id"
$els.Range("F3").Value  = "1"
$els.Range("J3").Value  = "string`n"
$els.Range("K3").Value  = "Unique id for inter-element referencing"
$els.Range("L3").Value  = "Unique id for the element within a resource (for internal references). This may be any string value that does not contain spaces."
$els.Range("AE3").Value = "Element.id"
$els.Range("AG3").Value = "1"
$els.Range("AJ3").Value = "n/a"

# Row 4: now describes "Extension.extension"
$els.Range("A4").Value  = "Extension.extension"
$els.Range("J4").Value  = "Extension`n"
$els.Range("L4").Value  = "An Extension"
$els.Range("AA4").Value = "value:url}`n"
$els.Range("AB4").Value = "Extensions are always sliced by (at least) url"
$els.Range("AD4").Value = "open"
$els.Range("AE4").Value = "Element.extension"
$els.Range("AI4").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`next-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"

# Row 5 ("Extension.url"): Min/Max style text reflows by one shared-string slot
$els.Range("E5").Value  = "1"
$els.Range("F5").Value  = "1"
$els.Range("AF5").Value = "1"
$els.Range("AG5").Value = "1"

# Row 6 ("Extension.value[x]"): same reflow
$els.Range("F6").Value  = "1"
$els.Range("J6").Value  = "string`n"
$els.Range("AG6").Value = "1"
